# Update column F (dSF) values for several rows in the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    4  = -5
    5  = 1
    6  = 4
    7  = 11
    8  = 2
    9  = -1
    10 = 0
    11 = -5
    12 = -2
    14 = -2
    21 = -2
    23 = -3
    24 = 5
    25 = -5
    27 = -3
    30 = -1
    32 = -6
    33 = 7
    35 = -1
    37 = -4
    38 = 0
    40 = -1
    42 = 3
    43 = 2
    44 = 6
    45 = -3
    46 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
